$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4225
$ws.Range("I19").Value = 1871.8334
$ws.Range("J19").Value = 6578.1665
$ws.Range("K19").Value = 1871.8334
$ws.Range("L19").Value = 6578.1665
$ws.Range("M19").Value = -1696.8334
$ws.Range("N19").Value = -6928.1665
$ws.Range("H86").Value = 5888.8887
$ws.Range("J86").Value = 6000
$ws.Range("L86").Value = 6000
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 5888.8887
$ws.Range("J89").Value = 6000
$ws.Range("L89").Value = 30000
$ws.Range("N89").Value = -41232
$ws.Range("H106").Value = 2411.875
$ws.Range("I106").Value = 2215.4614
$ws.Range("J106").Value = 3263
$ws.Range("K106").Value = 2215.4614
$ws.Range("L106").Value = 3263
$ws.Range("M106").Value = -1584.4614
$ws.Range("N106").Value = -4525
$ws.Range("H125").Value = 6000
$ws.Range("J125").Value = 6000
$ws.Range("L125").Value = 54000
$ws.Range("N125").Value = -58920
$ws.Range("H138").Value = 9013876
$ws.Range("I138").Value = 1170.4445
$ws.Range("K138").Value = 3511.3335
$ws.Range("M138").Value = 1628.6665

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10422207
$ws.Range("I32").Value = 13515731
$ws.Range("K32").Value = 13515731
$ws.Range("M32").Value = -13515444
$ws.Range("H45").Value = 1694.9
$ws.Range("I45").Value = 1266.5
$ws.Range("K45").Value = 1266.5
$ws.Range("M45").Value = -889.5
$ws.Range("H61").Value = 45457988
$ws.Range("I61").Value = 58825016
$ws.Range("J61").Value = 10102.8
$ws.Range("K61").Value = 58825016
$ws.Range("L61").Value = 10102.8
$ws.Range("M61").Value = -58824804
$ws.Range("N61").Value = -10526.8
$ws.Range("H63").Value = 4403.4116
$ws.Range("I63").Value = 2990.68
$ws.Range("J63").Value = 8327.666999999999
$ws.Range("K63").Value = 2990.68
$ws.Range("L63").Value = 8327.666999999999
$ws.Range("M63").Value = -2304.68
$ws.Range("N63").Value = -9699.666999999999
$ws.Range("H66").Value = 4403.4116
$ws.Range("I66").Value = 2990.68
$ws.Range("J66").Value = 8327.666999999999
$ws.Range("K66").Value = 14953.4
$ws.Range("L66").Value = 41638.335
$ws.Range("M66").Value = -11521.4
$ws.Range("N66").Value = -48502.335
$ws.Range("H97").Value = 1299.7142
$ws.Range("I97").Value = 384.30768
$ws.Range("J97").Value = 2787.25
$ws.Range("K97").Value = 384.30768
$ws.Range("L97").Value = 2787.25
$ws.Range("M97").Value = 111.69232
$ws.Range("N97").Value = -3779.25
$ws.Range("H102").Value = 3698.7
$ws.Range("I102").Value = 3122
$ws.Range("K102").Value = 3122
$ws.Range("M102").Value = -1500
$ws.Range("H103").Value = 59500
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344
$ws.Range("H132").Value = 29414020
$ws.Range("I132").Value = 2135.2258
$ws.Range("J132").Value = 333336830
$ws.Range("K132").Value = 6405.6774
$ws.Range("L132").Value = 1000010490
$ws.Range("M132").Value = -3875.6774
$ws.Range("N132").Value = -1000015550
$ws.Range("H136").Value = 45457988
$ws.Range("I136").Value = 58825016
$ws.Range("J136").Value = 10102.8
$ws.Range("K136").Value = 176475048
$ws.Range("L136").Value = 30308.4
$ws.Range("M136").Value = -176472498
$ws.Range("N136").Value = -35408.39999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 38288.43
$ws.Range("J81").Value = 38288.43
$ws.Range("L81").Value = 38288.43
$ws.Range("N81").Value = -40410.43
$ws.Range("H84").Value = 38288.43
$ws.Range("J84").Value = 38288.43
$ws.Range("L84").Value = 114865.29
$ws.Range("N84").Value = -125473.29
$ws.Range("H86").Value = 12655.033
$ws.Range("I86").Value = 6394.24
$ws.Range("J86").Value = 43959
$ws.Range("K86").Value = 6394.24
$ws.Range("L86").Value = 43959
$ws.Range("M86").Value = -5271.24
$ws.Range("N86").Value = -46205
$ws.Range("H89").Value = 12655.033
$ws.Range("I89").Value = 6394.24
$ws.Range("J89").Value = 43959
$ws.Range("K89").Value = 31971.2
$ws.Range("L89").Value = 219795
$ws.Range("M89").Value = -26355.2
$ws.Range("N89").Value = -231027
$ws.Range("H94").Value = 2096.96
$ws.Range("I94").Value = 1976.7858
$ws.Range("K94").Value = 1976.7858
$ws.Range("M94").Value = -1525.7858
$ws.Range("H105").Value = 25839.75
$ws.Range("I105").Value = 34003
$ws.Range("K105").Value = 34003
$ws.Range("M105").Value = -32256

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 37743
$ws.Range("J68").Value = 39838
$ws.Range("L68").Value = 39838
$ws.Range("N68").Value = -41336
$ws.Range("H71").Value = 37743
$ws.Range("J71").Value = 39838
$ws.Range("L71").Value = 119514
$ws.Range("N71").Value = -127002
$ws.Range("H86").Value = 5140.1665
$ws.Range("I86").Value = 5210.5
$ws.Range("K86").Value = 5210.5
$ws.Range("M86").Value = -4087.5
$ws.Range("H89").Value = 5140.1665
$ws.Range("I89").Value = 5210.5
$ws.Range("K89").Value = 26052.5
$ws.Range("M89").Value = -20436.5
$ws.Range("H93").Value = 14099.4
$ws.Range("I93").Value = 12888.223
$ws.Range("K93").Value = 12888.223
$ws.Range("M93").Value = -11016.223
$ws.Range("H122").Value = 1124.6552
$ws.Range("I122").Value = 1127.4
$ws.Range("K122").Value = 3382.2
$ws.Range("M122").Value = -932.2000000000003
$ws.Range("H134").Value = 1365.5834
$ws.Range("I134").Value = 1217
$ws.Range("K134").Value = 3651
$ws.Range("M134").Value = -1116

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 22480.455
$ws.Range("I56").Value = 22480.455
$ws.Range("K56").Value = 22480.455
$ws.Range("M56").Value = -21950.455
$ws.Range("H80").Value = 37999.332
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 37999.332
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 113997.996
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -115869.996
$ws.Range("H83").Value = 37999.332
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 37999.332
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 341993.988
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -351353.988
$ws.Range("H131").Value = 33766.57
$ws.Range("J131").Value = 5352.56
$ws.Range("L131").Value = 16057.68
$ws.Range("N131").Value = -26137.68
$ws.Range("H132").Value = 1711969.9
$ws.Range("J132").Value = 3706718.5
$ws.Range("L132").Value = 33360466.5
$ws.Range("N132").Value = -33365526.5
$ws.Range("H139").Value = 2101.842
$ws.Range("I139").Value = 2042.4
$ws.Range("J139").Value = 2324.75
$ws.Range("K139").Value = 6127.200000000001
$ws.Range("L139").Value = 6974.25
$ws.Range("M139").Value = -987.2000000000007
$ws.Range("N139").Value = -17254.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 16999.834
$ws.Range("I49").Value = 13000
$ws.Range("K49").Value = 13000
$ws.Range("M49").Value = -12816
$ws.Range("H80").Value = 2378.5454
$ws.Range("I80").Value = 2366.4
$ws.Range("K80").Value = 2366.4
$ws.Range("M80").Value = -1368.4
$ws.Range("H83").Value = 2378.5454
$ws.Range("I83").Value = 2366.4
$ws.Range("K83").Value = 11832
$ws.Range("M83").Value = -6840
$ws.Range("H97").Value = 1635.125
$ws.Range("I97").Value = 583.25
$ws.Range("J97").Value = 2687
$ws.Range("K97").Value = 583.25
$ws.Range("L97").Value = 2687
$ws.Range("M97").Value = -87.25
$ws.Range("N97").Value = -3679
$ws.Range("H113").Value = 3141.1365
$ws.Range("I113").Value = 1977.1538
$ws.Range("K113").Value = 1977.1538
$ws.Range("M113").Value = 192.8462

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5040.5454
$ws.Range("I61").Value = 4226.1665
$ws.Range("J61").Value = 6017.8
$ws.Range("K61").Value = 4226.1665
$ws.Range("L61").Value = 6017.8
$ws.Range("M61").Value = -4024.1665
$ws.Range("N61").Value = -6421.8
$ws.Range("H100").Value = 4075.25
$ws.Range("I100").Value = 2999
$ws.Range("K100").Value = 2999
$ws.Range("M100").Value = -2458
$ws.Range("H113").Value = 5040.5454
$ws.Range("I113").Value = 4226.1665
$ws.Range("J113").Value = 6017.8
$ws.Range("K113").Value = 4226.1665
$ws.Range("L113").Value = 6017.8
$ws.Range("M113").Value = -2056.1665
$ws.Range("N113").Value = -10357.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 20615
$ws.Range("I52").Value = 15645
$ws.Range("J52").Value = 40495
$ws.Range("K52").Value = 15645
$ws.Range("L52").Value = 40495
$ws.Range("M52").Value = -15419
$ws.Range("N52").Value = -40947
$ws.Range("H96").Value = 7217.625
$ws.Range("J96").Value = 9968.333000000001
$ws.Range("L96").Value = 9968.333000000001
$ws.Range("N96").Value = -12714.333
$ws.Range("H122").Value = 41667852
$ws.Range("I122").Value = 47620140
$ws.Range("K122").Value = 142860420
$ws.Range("M122").Value = -142857970
$ws.Range("H136").Value = 1287.921
$ws.Range("I136").Value = 1058.6285
$ws.Range("K136").Value = 3175.8855
$ws.Range("M136").Value = -625.8855000000003
